# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-03 11:14:41
#
# Applies the updated attendance-recording data to the
# "Session Analysis Results" sheet: refreshed "Recorded By" lists
# (re-ordered), updated Class/Group statistics numbers, and the
# B4/PHYSIOLOGY/session-1 row that moved from "Not Recorded" to "Recorded".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextValue {
    # NOTE: this runtime's PowerShell engine does not bind named (-Param)
    # arguments correctly, so this function (and every other call in this
    # script) relies on POSITIONAL parameters only:
    #   $args[0] = cell reference (e.g. "L9")
    #   $args[1] = text to store in the cell
    #   $args[2] = reference cell to copy the original number format/style from
    param([string]$cellRef, [string]$text, [string]$templateRef)
    # Force the cell to remain a text value (e.g. "14.2%") instead of being
    # auto-coerced into a numeric/percentage value, then restore the
    # original cell style (fill/font/alignment) by pasting formats from an
    # existing cell that already carries that exact style.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($templateRef).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# G2 / G24 - ANATOMY A1 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g2 = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G2").Value = $g2
$ws.Range("G24").Value = $g2

# ---------------------------------------------------------------------
# Class Statistics block (K5:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 13
$ws.Range("L8").Value = 138
Set-TextValue "L9" "14.2%" "L5"
Set-TextValue "L10" "29.5%" "L5"

# ---------------------------------------------------------------------
# Group Statistics block (K14:S22)
# ---------------------------------------------------------------------
# Year 2 / A1 (row 15)
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 17

# Year 2 / A2 (row 16)
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 17

# Year 2 / A3 (row 17)
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 18

# Year 2 / A4 (row 18)
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 18

# Year 2 / B4 (row 22)
$ws.Range("O22").Value = 4
$ws.Range("P22").Value = 0
Set-TextValue "R22" "18.2%" "L5"
Set-TextValue "S22" "11.9%" "L5"

# ---------------------------------------------------------------------
# G18 / G40 - PHYSIOLOGY A1 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g18 = "shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G18").Value = $g18
$ws.Range("G40").Value = $g18

# ---------------------------------------------------------------------
# Row 19: Year 2 / A1 / PHYSIOLOGY / session 2 -> now "Not Recorded"
# (style changes from "Pending" (s=4) to "Not Recorded" (s=6); copy the
# format from row 7, which already carries the "Not Recorded" style.)
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A19:I19").PasteSpecial($xlPasteFormats)
$ws.Range("I19").Value = "Not Recorded"

# ---------------------------------------------------------------------
# Row 36: Year 2 / A2 / PHARMACOLOGY / session 1 -> now "Not Recorded"
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A36:I36").PasteSpecial($xlPasteFormats)
$ws.Range("I36").Value = "Not Recorded"

# ---------------------------------------------------------------------
# G52 / G74 - HISTOLOGY A3 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g52 = "mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G52").Value = $g52
$ws.Range("G74").Value = $g52

# ---------------------------------------------------------------------
# Row 54: Year 2 / A3 / MICROBIOLOGY / session 1 -> now "Not Recorded"
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A54:I54").PasteSpecial($xlPasteFormats)
$ws.Range("I54").Value = "Not Recorded"

# ---------------------------------------------------------------------
# G62 / G84 - PHYSIOLOGY A3 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g62 = "wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G62").Value = $g62
$ws.Range("G84").Value = $g62

# ---------------------------------------------------------------------
# Row 73: Year 2 / A4 / CARDIOLOGY / session 1 -> now "Not Recorded"
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A73:I73").PasteSpecial($xlPasteFormats)
$ws.Range("I73").Value = "Not Recorded"

# ---------------------------------------------------------------------
# G96 / G118 - HISTOLOGY B1/B2 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g96 = "Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G96").Value = $g96
$ws.Range("G118").Value = $g96

# ---------------------------------------------------------------------
# G98 / G120 - MICROBIOLOGY B1/B2 "Recorded By" list reordered
# ---------------------------------------------------------------------
$g98 = "yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G98").Value = $g98
$ws.Range("G120").Value = $g98

# ---------------------------------------------------------------------
# G134 - ANATOMY B3 "Recorded By" list reordered
# ---------------------------------------------------------------------
$ws.Range("G134").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ---------------------------------------------------------------------
# G150 - PHYSIOLOGY B3 "Recorded By" list reordered + new recorder added
# ---------------------------------------------------------------------
$ws.Range("G150").Value = "nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"

# ---------------------------------------------------------------------
# G156 - ANATOMY B4 "Recorded By" list reordered
# ---------------------------------------------------------------------
$ws.Range("G156").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 172: Year 2 / B4 / PHYSIOLOGY / session 1 -> now "Recorded"
# (style changes from "Not Recorded" (s=6) to "Recorded" (s=2); copy the
# format from row 2, which already carries the "Recorded" style.)
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A172:I172").PasteSpecial($xlPasteFormats)
$ws.Range("G172").Value = "nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("H172").Value = "6/226"
$ws.Range("I172").Value = "Recorded"

$excel.CutCopyMode = $false
Write-Host "Edit complete."
